# Update the "Time:" timestamp embedded in the OLS regression summary text
# (cell B2) on every worksheet, from 20:51:38 to 20:59:40.
$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = $cell.Value2
    if ($null -ne $text -and $text.Contains("20:51:38")) {
        $cell.Value2 = $text.Replace("20:51:38", "20:59:40")
    }
}
